$d = $word.ActiveDocument

# The last paragraph in the body is the "Google Authentication." list item.
# Insert a brand-new list item right after it, in the same list
# (pStyle "ListParagraph", ilvl 0, numId 1) with the new bullet text.
$lastPara = $d.Paragraphs.Last
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$newParaXml = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Used iframe to embed the preview of a movie</w:t></w:r></w:p>'

$insertPoint.InsertXML($newParaXml)
